# Add files via upload
# Fixes the Image column references on the "7 Inch" pizza rows so that each
# dish points at its correct (renamed/re-extensioned) image file, and
# corrects the "Pizza Corn 10'.png" shared string (drop stray apostrophe).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rows 10-13: re-point each 7" pizza item to its corrected image filename.
# (Order matters for shared-string table layout - matches original author edit order.)
$ws.Range("D12").Value = "7 Inch Pizza Panner.jpg"
$ws.Range("D11").Value = "7 Inch Pizza Onion and Capsicum.jpg"
$ws.Range("D6").Value = "Pizza Corn 10.png"
$ws.Range("D10").Value = "7 Inch Pizza Margarita.png"
$ws.Range("D13").Value = "7 Inch Pizza Corn.png"

# Update the remembered cell selection on Sheet1.
$ws.Range("H11").Select()
